$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Solitario"
$ws.Range("B12").Value = "Ulises el amor de mi vida"
$ws.Range("C12").Value = 732
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = "2025-11-11 00:29:36"

$ws.Range("A13").Value = "Multijugador"
$ws.Range("B13").Value = "Guillermina"
$ws.Range("C13").Value = 283
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = "2025-11-11 00:30:29"
